$d = $word.ActiveDocument

$d.Content.Find.Execute("литературных", $true, $false, $false, $false, $false,
                         $true, 1, $false, "использованных", 2)
